# Apply the commit's change: the long "questions = [...]" text (formerly in
# A2, stored as a shared string) is pretty-printed (json.dumps(questions, indent=4))
# and moved into A1, replacing the old placeholder value (0) that used a bold /
# bordered / centered style. A2 is removed entirely, and A1 reverts to the
# default (unstyled) cell formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lines = @(
    "questions = [",
    "    {",
    "        `"title`": `"True or false: Spark supports caching data in memory for individual clusters.`",",
    "        `"ques_type`": 11,",
    "        `"options`": [",
    "            `"true`",",
    "            `"false`"",
    "        ],",
    "        `"score`": `"True`"",
    "    },",
    "    {",
    "        `"title`": `"What will be the output of the code below?spark.createDataFrame([('Spark', 1)], ['Course', 'Rank']).collect()`",",
    "        `"ques_type`": 2,",
    "        `"options`": [",
    "            `"[Row(Course='Spark', Rank=1)]`",",
    "            `"[Row(_1='Spark', _2=1)]`",",
    "            `"[Row(Spark=Course, 1=Rank)]`",",
    "            `"[Row(Spark=_1, 1=_2)]`"",
    "        ],",
    "        `"score`": `"[Row(Course='Spark', Rank=1)]`"",
    "    },",
    "    {",
    "        `"title`": `"When you ran the following code to create two DataFrames named \u201cdf1\u201d and \u201cdf2,\u201d it resulted in an analysis error.  df1 = spark.createDataFrame([[1, 2]], [\`"col1\`", \`"col2\`"])df2 = spark.createDataFrame([[2, 3]], [\`"col2\`", \`"col3\`"]) You expected the output as per Table A, shown below. Which operation should you perform to obtain the output you expected? Table A`",",
    "        `"ques_type`": 2,",
    "        `"options`": [",
    "            `"df1.unionByName(df2).show()`",",
    "            `"df1.join(df2).show()`",",
    "            `"df2.unionByName(df1).show()`",",
    "            `"df1.unionByName(df2, allowMissingColumns=True).show()`",",
    "            `"df2.join(df1).show()`"",
    "        ],",
    "        `"score`": `"df1.unionByName(df2, allowMissingColumns=True).show()`"",
    "    },",
    "    {",
    "        `"title`": `"What will be the output of the code below?sparkContext.parallelize([(1, 'apple'), (1, 'mango')]).countByKey().items()`",",
    "        `"ques_type`": 2,",
    "        `"options`": [",
    "            `"[((1, 'apple'), 1), ((2, 'mango'), 1)]`",",
    "            `"[(1, 2)]`",",
    "            `"[(1, 'apple'), (1, 'mango')]`",",
    "            `"[(\u2018apple\u2019, \u2018mango\u2019)]`"",
    "        ],",
    "        `"score`": `"[(1, 2)]`"",
    "    }",
    "]"
)
$newText = $lines -join "`n"

# Put the new pretty-printed text into A1.
$ws.Range("A1").Value2 = $newText

# A1 previously carried a bold/bordered/centered style; reset it back to the
# workbook's default "Normal" style so no custom formatting remains.
$ws.Range("A1").Style = "Normal"

# Setting a long, multi-line value auto-expands the row height; restore it to
# fit the (now unstyled, unwrapped) single-line default presentation.
$ws.Rows(1).AutoFit()

# The old A2 cell (which held the shared string) is no longer needed - remove
# the whole row so the sheet shrinks back down to just A1.
$ws.Rows(2).Delete()
